$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '29.394.33'
$ws.Range("E2").Value = '  -1.98%  '

Set-TextValue $ws.Range("D3") '1.995.44'
$ws.Range("E3").Value = '  -5.66%  '

$ws.Range("E4").Value = '  -0.01%  '

Set-TextValue $ws.Range("D5") '331.30'
$ws.Range("E5").Value = '  -4.27%  '

$ws.Range("E6").Value = '  -0.05%  '

Set-TextValue $ws.Range("D7") '0.4925'
$ws.Range("E7").Value = '  -5.13%  '

Set-TextValue $ws.Range("D8") '0.4177'
$ws.Range("E8").Value = '  -6.06%  '

Set-TextValue $ws.Range("D9") '53.31'
$ws.Range("E9").Value = '  -0.93%  '

Set-TextValue $ws.Range("D10") '0.08778'
$ws.Range("E10").Value = '  -6.37%  '

Set-TextValue $ws.Range("D11") '1.110'
$ws.Range("E11").Value = '  -6.08%  '

Set-TextValue $ws.Range("D12") '2.133.56'
$ws.Range("E12").Value = '  +0.99%  '

Set-TextValue $ws.Range("D13") '23.03'
$ws.Range("E13").Value = '  -8.72%  '

Set-TextValue $ws.Range("D14") '8.020'
$ws.Range("E14").Value = '  -6.37%  '

Set-TextValue $ws.Range("D15") '6.463'
$ws.Range("E15").Value = '  -6.45%  '

Set-TextValue $ws.Range("D16") '96.02'
$ws.Range("E16").Value = '  -6.69%  '

$ws.Range("E17").Value = '  +0.07%  '

$ws.Range("E18").Value = '  -5.03%  '

Set-TextValue $ws.Range("D19") '0.06628'
$ws.Range("E19").Value = '  -1.02%  '

Set-TextValue $ws.Range("D20") '19.43'
$ws.Range("E20").Value = '  -9.58%  '

$ws.Range("E21").Value = '  +0.06%  '

Set-TextValue $ws.Range("D22") '5.951'
$ws.Range("E22").Value = '  -5.52%  '

Set-TextValue $ws.Range("D23") '29.441.13'
$ws.Range("E23").Value = '  -1.89%  '

Set-TextValue $ws.Range("D24") '11.77'
$ws.Range("E24").Value = '  -7.44%  '

Set-TextValue $ws.Range("D25") '2.283'
$ws.Range("E25").Value = '  -1.56%  '

Set-TextValue $ws.Range("D26") '2.307.47'
$ws.Range("E26").Value = '  -2.28%  '

Set-TextValue $ws.Range("D27") '6.655'
$ws.Range("E27").Value = '  +0.40%  '

Set-TextValue $ws.Range("D28") '157.28'
$ws.Range("E28").Value = '  -3.27%  '

Set-TextValue $ws.Range("D29") '20.46'
$ws.Range("E29").Value = '  -7.42%  '

Set-TextValue $ws.Range("D30") '2.340'
$ws.Range("E30").Value = '  -7.63%  '

Set-TextValue $ws.Range("D31") '126.48'
$ws.Range("E31").Value = '  -5.72%  '

Set-TextValue $ws.Range("D32") '1.047'
$ws.Range("E32").Value = '  -9.06%  '

Set-TextValue $ws.Range("D33") '0.09903'
$ws.Range("E33").Value = '  -6.23%  '

Set-TextValue $ws.Range("D34") '1.547'
$ws.Range("E34").Value = '  -13.43%  '

Set-TextValue $ws.Range("D35") '5.804'
$ws.Range("E35").Value = '  -7.14%  '

Set-TextValue $ws.Range("D36") '3.778'
$ws.Range("E36").Value = '  -4.79%  '

Set-TextValue $ws.Range("D37") '9.578'
$ws.Range("E37").Value = '  -11.13%  '

Set-TextValue $ws.Range("D38") '0.02442'
$ws.Range("E38").Value = '  -6.50%  '

Set-TextValue $ws.Range("D39") '0.06351'
$ws.Range("E39").Value = '  -7.36%  '

Set-TextValue $ws.Range("D40") '1.279'
$ws.Range("E40").Value = '  -4.13%  '

$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D41") '11.70'
$ws.Range("E41").Value = '  -7.88%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D42") '0.6475'
$ws.Range("E42").Value = '  -8.51%  '

Set-TextValue $ws.Range("D43") '0.2058'
$ws.Range("E43").Value = '  -8.00%  '

Set-TextValue $ws.Range("D44") '1.007'
$ws.Range("E44").Value = '  +0.01%  '

Set-TextValue $ws.Range("D45") '0.6301'
$ws.Range("E45").Value = '  -7.80%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D46") '13.31'
$ws.Range("E46").Value = '  -9.56%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D47") '2.187'
$ws.Range("E47").Value = '  -7.78%  '

Set-TextValue $ws.Range("D48") '1.260'
$ws.Range("E48").Value = '  -0.76%  '

Set-TextValue $ws.Range("D49") '3.541'
$ws.Range("E49").Value = '  -2.37%  '

Set-TextValue $ws.Range("D50") '0.06984'
$ws.Range("E50").Value = '  -1.77%  '

Set-TextValue $ws.Range("D51") '1.142'
$ws.Range("E51").Value = '  -3.40%  '
